$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Preserve original cell style while forcing the written value to stay
    # a text cell (matches source data where Price/Volume are inline strings,
    # even when the text looks numeric, e.g. trailing-zero prices "71.10").
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "36.236.12"
$ws.Range("E2").Value = "  -4.02%  "

Set-TextValue $ws.Range("D3") "1.966.95"
$ws.Range("E3").Value = "  -4.19%  "

$ws.Range("E4").Value = "  +0.22%  "

Set-TextValue $ws.Range("D5") "242.09"
$ws.Range("E5").Value = "  -4.27%  "

Set-TextValue $ws.Range("D6") "0.624"
$ws.Range("E6").Value = "  -4.35%  "

Set-TextValue $ws.Range("D7") "59.72"
$ws.Range("E7").Value = "  -9.03%  "

$ws.Range("E8").Value = "  +0.04%  "

Set-TextValue $ws.Range("D9") "0.373"
$ws.Range("E9").Value = "  -1.61%  "

Set-TextValue $ws.Range("D10") "56.96"
$ws.Range("E10").Value = "  -4.86%  "

Set-TextValue $ws.Range("D11") "0.0802"
$ws.Range("E11").Value = "  +5.21%  "

$ws.Range("E12").Value = "  -0.52%  "

Set-TextValue $ws.Range("D13") "0.858"
$ws.Range("E13").Value = "  -7.45%  "

Set-TextValue $ws.Range("D14") "22.46"
$ws.Range("E14").Value = "  +9.57%  "

Set-TextValue $ws.Range("D15") "14.02"
$ws.Range("E15").Value = "  -7.76%  "

Set-TextValue $ws.Range("D16") "2.253.03"
$ws.Range("E16").Value = "  -4.07%  "

$ws.Range("E17").Value = "  -2.94%  "

Set-TextValue $ws.Range("D18") "1.967.09"
$ws.Range("E18").Value = "  -3.89%  "

Set-TextValue $ws.Range("D19") "36.102.67"
$ws.Range("E19").Value = "  -4.03%  "

Set-TextValue $ws.Range("D20") "71.10"
$ws.Range("E20").Value = "  -3.72%  "

Set-TextValue $ws.Range("D21") "0.0₃0858"
$ws.Range("E21").Value = "  -2.67%  "

Set-TextValue $ws.Range("D22") "237.30"
$ws.Range("E22").Value = "  -0.38%  "

Set-TextValue $ws.Range("D23") "5.21"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("E24").Value = "  +0.11%  "

Set-TextValue $ws.Range("D25") "2.52"
$ws.Range("E25").Value = "  -5.88%  "

$ws.Range("E26").Value = "  -4.53%  "

$ws.Range("E27").Value = "  +1.70%  "

Set-TextValue $ws.Range("D28") "160.44"
$ws.Range("E28").Value = "  +0.09%  "

Set-TextValue $ws.Range("D29") "19.85"
$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("E30").Value = "  +11.71%  "

$ws.Range("E31").Value = "  -2.36%  "

$ws.Range("E32").Value = "  -7.41%  "

$ws.Range("E33").Value = "  -5.81%  "

Set-TextValue $ws.Range("D34") "0.0620"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("E35").Value = "  -7.50%  "

$ws.Range("E36").Value = "  +5.38%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D38") "2.29"
$ws.Range("E38").Value = "  -6.32%  "

Set-TextValue $ws.Range("D39") "1.81"
$ws.Range("E39").Value = "  -2.46%  "

Set-TextValue $ws.Range("D40") "3.10"
$ws.Range("E40").Value = "  +9.27%  "

Set-TextValue $ws.Range("D41") "0.0987"
$ws.Range("E41").Value = "  -3.88%  "

$ws.Range("E42").Value = "  -1.06%  "

Set-TextValue $ws.Range("D43") "2.86"
$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("E44").Value = "  -3.17%  "

$ws.Range("E45").Value = "  -4.79%  "

Set-TextValue $ws.Range("D46") "92.27"
$ws.Range("E46").Value = "  -3.61%  "

Set-TextValue $ws.Range("D47") "15.99"
$ws.Range("E47").Value = "  -6.19%  "

$ws.Range("E48").Value = "  -7.16%  "

Set-TextValue $ws.Range("D49") "1.337.22"
$ws.Range("E49").Value = "  -6.51%  "

$ws.Range("E50").Value = "  -3.72%  "

Set-TextValue $ws.Range("D51") "2.147.13"
$ws.Range("E51").Value = "  -3.97%  "
